# Commit: "Created Function for Gaussian Quadrature Scheme, and exported it
#          to the Averaged Intensities files."
#
# 1. Rename the worksheet (and its tab name) from "UniformA-HW20.xpc" to "UniformA".
# 2. Append a new data row (row 16) that mirrors the existing rows' layout:
#    A16 = 14 (index), B16 = "HexGrid-60degTilt5degRes" (same label text as row 15),
#    C16:P16 = 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet
$ws.Name = "UniformA"

# Seed row 16 with row 15's formatting (border/alignment/font) by copying the
# whole row, then overwrite with the new row's values.
$ws.Range("A15:P15").Copy($ws.Range("A16:P16"))

$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "HexGrid-60degTilt5degRes"
for ($col = 3; $col -le 16; $col++) {
    $ws.Cells.Item(16, $col).Value = 1
}
